$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.000", "0.00001092") are preserved exactly as text,
# matching the original inlineStr cell types.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '21.660.68'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '1.534.67'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '288.42'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = '0.3923'
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("D8").Value = '0.3171'
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").Value = '42.46'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").Value = '0.07177'
$ws.Range("E10").Value = '  -2.59%  '
$ws.Range("D11").Value = '1.053'
$ws.Range("E11").Value = '  -6.30%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '5.671'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").Value = '18.59'
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.560.57'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '6.581'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").Value = '0.00001092'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").Value = '0.06589'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '83.78'
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '6.120'
$ws.Range("E21").Value = '  -4.27%  '
$ws.Range("D22").Value = '15.47'
$ws.Range("E22").Value = '  -3.28%  '
$ws.Range("D23").Value = '10.73'
$ws.Range("E23").Value = '  -6.51%  '
$ws.Range("D24").Value = '2.352'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").Value = '21.654.90'
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("D26").Value = '2.357'
$ws.Range("E26").Value = '  -7.95%  '
$ws.Range("D27").Value = '149.43'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = '18.31'
$ws.Range("E28").Value = '  -3.00%  '
$ws.Range("D29").Value = '4.844'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = '1.731.79'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '117.03'
$ws.Range("E31").Value = '  -3.20%  '
$ws.Range("D32").Value = '6.047'
$ws.Range("E32").Value = '  +3.27%  '
$ws.Range("D33").Value = '0.9394'
$ws.Range("E33").Value = '  -15.53%  '
$ws.Range("D34").Value = '0.08151'
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("D35").Value = '8.514'
$ws.Range("E35").Value = '  -8.28%  '
$ws.Range("D36").Value = '5.160'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").Value = '0.06024'
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("D38").Value = '0.02221'
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("D39").Value = '1.455'
$ws.Range("E39").Value = '  -14.65%  '
$ws.Range("D40").Value = '0.2021'
$ws.Range("E40").Value = '  -4.05%  '
$ws.Range("D41").Value = '1.177'
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("D42").Value = '10.96'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '0.5761'
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("D45").Value = '12.99'
$ws.Range("E45").Value = '  -4.33%  '
$ws.Range("D46").Value = '3.708'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").Value = '0.5501'
$ws.Range("E47").Value = '  -4.40%  '
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = '1.879'
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("D50").Value = '116.23'
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("D51").Value = '0.06691'
$ws.Range("E51").Value = '  -2.96%  '
